$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value2 = 1338.6
$ws.Range("I94").Value2 = 1338.6
$ws.Range("K94").Value2 = 1338.6
$ws.Range("M94").Value2 = -887.5999999999999

$ws.Range("H135").Value2 = 511.63333
$ws.Range("I135").Value2 = 353.52173
$ws.Range("J135").Value2 = 1031.1428
$ws.Range("K135").Value2 = 3181.69557
$ws.Range("L135").Value2 = 9280.2852
$ws.Range("M135").Value2 = -646.6955699999999
$ws.Range("N135").Value2 = -14350.2852

$ws.Range("H138").Value2 = 3997.9854
$ws.Range("I138").Value2 = 5368.28
$ws.Range("J138").Value2 = 3201.3022
$ws.Range("K138").Value2 = 16104.84
$ws.Range("L138").Value2 = 9603.9066
$ws.Range("M138").Value2 = -10964.84
$ws.Range("N138").Value2 = -19883.9066

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value2 = 185800
$ws.Range("I34").Value2 = 219750
$ws.Range("J34").Value2 = 50000
$ws.Range("K34").Value2 = 219750
$ws.Range("L34").Value2 = 50000
$ws.Range("M34").Value2 = -219479
$ws.Range("N34").Value2 = -50542

$ws.Range("H52").Value2 = 39998

$ws.Range("H88").Value2 = 1456.0834
$ws.Range("J88").Value2 = 1202
$ws.Range("L88").Value2 = 1202
$ws.Range("N88").Value2 = -2014

$ws.Range("H91").Value2 = 1456.0834
$ws.Range("J91").Value2 = 1202
$ws.Range("L91").Value2 = 1202
$ws.Range("N91").Value2 = -4010

$ws.Range("H97").Value2 = 6617.5713
$ws.Range("J97").Value2 = 3233.4285
$ws.Range("L97").Value2 = 3233.4285
$ws.Range("N97").Value2 = -4225.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 2412.9412
$ws.Range("I86").Value2 = 1674.3
$ws.Range("J86").Value2 = 3468.1428
$ws.Range("K86").Value2 = 1674.3
$ws.Range("L86").Value2 = 3468.1428
$ws.Range("M86").Value2 = -551.3
$ws.Range("N86").Value2 = -5714.1428

$ws.Range("H89").Value2 = 2412.9412
$ws.Range("I89").Value2 = 1674.3
$ws.Range("J89").Value2 = 3468.1428
$ws.Range("K89").Value2 = 8371.5
$ws.Range("L89").Value2 = 17340.714
$ws.Range("M89").Value2 = -2755.5
$ws.Range("N89").Value2 = -28572.714

$ws.Range("H94").Value2 = 5223.4287
$ws.Range("I94").Value2 = 2910.8
$ws.Range("K94").Value2 = 2910.8
$ws.Range("M94").Value2 = -2459.8

$ws.Range("H134").Value2 = 1901
$ws.Range("I134").Value2 = 1733.5172
$ws.Range("J134").Value2 = 2508.125
$ws.Range("K134").Value2 = 5200.5516
$ws.Range("L134").Value2 = 7524.375
$ws.Range("M134").Value2 = -2665.5516
$ws.Range("N134").Value2 = -12594.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 2262
$ws.Range("J58").Value2 = 1972.2727
$ws.Range("L58").Value2 = 1972.2727
$ws.Range("N58").Value2 = -2378.2727

$ws.Range("H134").Value2 = 1568.76
$ws.Range("I134").Value2 = 1079
$ws.Range("K134").Value2 = 3237
$ws.Range("M134").Value2 = -702

$ws.Range("H136").Value2 = 2262
$ws.Range("J136").Value2 = 1972.2727
$ws.Range("L136").Value2 = 5916.8181
$ws.Range("N136").Value2 = -11016.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value2 = 250015.75
$ws.Range("I23").Value2 = 26.5
$ws.Range("J23").Value2 = 500005
$ws.Range("K23").Value2 = 79.5
$ws.Range("L23").Value2 = 1500015
$ws.Range("M23").Value2 = 155.5
$ws.Range("N23").Value2 = -1500485

$ws.Range("H50").Value2 = 7000
$ws.Range("I50").Value2 = 0
$ws.Range("J50").Value2 = 7000
$ws.Range("K50").Value2 = 0
$ws.Range("M50").Value2 = 21000
$ws.Range("N50").Value2 = -21962

$ws.Range("H53").Value2 = 7000
$ws.Range("I53").Value2 = 0
$ws.Range("J53").Value2 = 7000
$ws.Range("K53").Value2 = 0
$ws.Range("M53").Value2 = 21000
$ws.Range("N53").Value2 = -21962

$ws.Range("H56").Value2 = 6911.875
$ws.Range("I56").Value2 = 6911.875
$ws.Range("K56").Value2 = 6911.875
$ws.Range("M56").Value2 = -6381.875

$ws.Range("H113").Value2 = 251
$ws.Range("I113").Value2 = 303.6
$ws.Range("J113").Value2 = 237.85
$ws.Range("K113").Value2 = 910.8000000000001
$ws.Range("L113").Value2 = 713.55
$ws.Range("M113").Value2 = 1259.2
$ws.Range("N113").Value2 = -5053.55

$ws.Range("H122").Value2 = 7105617.5
$ws.Range("J122").Value2 = 1701763.6
$ws.Range("L122").Value2 = 15315872.4
$ws.Range("N122").Value2 = -15320772.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value2 = 52777
$ws.Range("J34").Value2 = 52777
$ws.Range("L34").Value2 = 52777
$ws.Range("N34").Value2 = -53313

$ws.Range("H43").Value2 = 3692.8572
$ws.Range("J43").Value2 = 4750
$ws.Range("L43").Value2 = 4750
$ws.Range("N43").Value2 = -5052

$ws.Range("H54").Value2 = 3000
$ws.Range("I54").Value2 = 3000
$ws.Range("K54").Value2 = 3000
$ws.Range("M54").Value2 = -2610

$ws.Range("H68").Value2 = 50000
$ws.Range("J68").Value2 = 50000
$ws.Range("L68").Value2 = 50000
$ws.Range("N68").Value2 = -51622

$ws.Range("H71").Value2 = 50000
$ws.Range("J71").Value2 = 50000
$ws.Range("L71").Value2 = 150000
$ws.Range("N71").Value2 = -158112

$ws.Range("H76").Value2 = 52777
$ws.Range("J76").Value2 = 52777
$ws.Range("L76").Value2 = 52777
$ws.Range("N76").Value2 = -53407

$ws.Range("H79").Value2 = 52777
$ws.Range("J79").Value2 = 52777
$ws.Range("L79").Value2 = 52777
$ws.Range("N79").Value2 = -54961

$ws.Range("H126").Value2 = 4794.5454
$ws.Range("I126").Value2 = 4499.5
$ws.Range("J126").Value2 = 4963.143
$ws.Range("K126").Value2 = 13498.5
$ws.Range("L126").Value2 = 14889.429
$ws.Range("M126").Value2 = -11028.5
$ws.Range("N126").Value2 = -19829.429

$ws.Range("H132").Value2 = 4100057.8
$ws.Range("I132").Value2 = 1718.3871
$ws.Range("J132").Value2 = 8335008.5
$ws.Range("K132").Value2 = 5155.1613
$ws.Range("L132").Value2 = 25005025.5
$ws.Range("M132").Value2 = -2625.1613
$ws.Range("N132").Value2 = -25010085.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value2 = 14999.5
$ws.Range("J33").Value2 = 10000
$ws.Range("L33").Value2 = 10000
$ws.Range("N33").Value2 = -10580

$ws.Range("H44").Value2 = 13998
$ws.Range("J44").Value2 = 13998
$ws.Range("L44").Value2 = 13998
$ws.Range("N44").Value2 = -14910

$ws.Range("H82").Value2 = 2276.7856
$ws.Range("J82").Value2 = 2798.2856
$ws.Range("L82").Value2 = 2798.2856
$ws.Range("N82").Value2 = -3520.2856

$ws.Range("H85").Value2 = 2276.7856
$ws.Range("J85").Value2 = 2798.2856
$ws.Range("L85").Value2 = 2798.2856
$ws.Range("N85").Value2 = -5294.2856

$ws.Range("H136").Value2 = 26820.45
$ws.Range("J136").Value2 = 1240.3704
$ws.Range("L136").Value2 = 3721.1112
$ws.Range("N136").Value2 = -8821.111199999999
